$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaxCalculation")

# A3:A6 were empty numeric cells; the target state stores literal text
# values "1".."4" in them (rule-index labels), not numbers and not
# formulas. A leading apostrophe is the standard Excel way to force a
# numeric-looking entry to be stored as text.
$ws.Range("A3").Value = "'1"
$ws.Range("A4").Value = "'2"
$ws.Range("A5").Value = "'3"
$ws.Range("A6").Value = "'4"
